$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "Fecha" (date) values between the two week-1 rows (2-3) and
# the two week-2 rows (6-7): 2021-12-29 <-> 2022-01-13
$ws.Range("D2").Value = 44574
$ws.Range("D3").Value = 44574
$ws.Range("D6").Value = 44559
$ws.Range("D7").Value = 44559
